$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.718.40'
$ws.Range("D3").Value = '1.739.39'
$ws.Range("E3").Value = '  -5.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.89'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -7.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4908'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -6.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.75'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -7.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2420'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -23.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06020'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -11.45%  '
$ws.Range("D11").Value = '1.733.63'
$ws.Range("E11").Value = '  -5.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06725'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -13.27%  '
$ws.Range("E13").Value = '  -20.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5936'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -23.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -12.63%  '
$ws.Range("E16").Value = '  -12.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '25.737.05'
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -16.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006354'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -19.65%  '
$ws.Range("D22").Value = '1.958.55'
$ws.Range("E22").Value = '  -5.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.922'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -14.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.131'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -13.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.861'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -15.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '135.68'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.846'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -16.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.419'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -15.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.35'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -15.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '101.08'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08161'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.688'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -11.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.334'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -17.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04370'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -10.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9993'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.668'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.032'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -9.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6066'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -16.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.762'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -10.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.069'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.70%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.08'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.95%  '
$ws.Range("E43").Value = '  -13.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7918'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -11.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3805'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -20.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.120'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -13.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.071'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -20.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05085'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -12.56%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.77'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -14.21%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.03'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -12.92%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.239'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -12.03%  '
